# Insert a new date column ("04dec2025") immediately before the existing
# "26nov2025" column (column E) on both worksheets, shifting the old
# E:I columns to F:J, then populate the new column with the day's counts.
#
# "crosstab" holds real numbers; "annot" mirrors the same numbers but
# stored as text, with 0 rendered as a blank cell (matching the existing
# convention already used for every other column on that sheet).

$wb = $excel.ActiveWorkbook

# New column values (crosstab / numeric), keyed by row number.
$newValues = @{
    2  = 1
    3  = 4
    4  = 0
    5  = 2
    6  = 2
    7  = 3
    8  = 1
    9  = 2
    10 = 3
    11 = 3
    12 = 7
    13 = 3
    14 = 4
    15 = 3
    16 = 1
    17 = 2
    18 = 6
    19 = 3
    20 = 3
    21 = 0
    22 = 3
    23 = 3
}

for ($i = 1; $i -le $wb.Worksheets.Count(); $i++) {
    $ws = $wb.Worksheets.Item($i)
    $name = $ws.Name()

    # Insert a new column at E; everything from E onward shifts right by one.
    $ws.Columns("E:E").Insert()

    # Header for the newly inserted column.
    $ws.Range("E1").Value = "04dec2025"

    if ($name -eq "annot") {
        # Text mirror sheet: store the value as text (0 -> blank), matching
        # the style already used by every other data column on this sheet.
        $col = $ws.Range("E2:E23")
        $col.NumberFormat = "@"
        foreach ($r in $newValues.Keys) {
            $v = $newValues[$r]
            if ($v -eq 0) {
                $ws.Cells.Item($r, 5).Value = ""
            } else {
                $ws.Cells.Item($r, 5).Value = [string]$v
            }
        }
        $col.NumberFormat = "General"
    } else {
        # Numeric crosstab sheet.
        foreach ($r in $newValues.Keys) {
            $ws.Cells.Item($r, 5).Value = $newValues[$r]
        }
    }
}
